$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$updates = @{
    "D2" = "24.415.02";
    "E2" = "  -2.03%  ";
    "D3" = "1.647.69";
    "E3" = "  -3.78%  ";
    "D4" = "1.003";
    "E4" = "  -0.09%  ";
    "D5" = "310.77";
    "E5" = "  -0.86%  ";
    "D6" = "0.9985";
    "E6" = "  +0.11%  ";
    "D7" = "0.3646";
    "E7" = "  -3.09%  ";
    "D8" = "46.74";
    "E8" = "  -5.99%  ";
    "D9" = "0.3245";
    "E9" = "  -6.58%  ";
    "D10" = "1.122";
    "E10" = "  -7.94%  ";
    "D11" = "0.07024";
    "E11" = "  -7.64%  ";
    "D12" = "0.9996";
    "E12" = "  -0.01%  ";
    "D13" = "5.956";
    "E13" = "  -6.44%  ";
    "D14" = "19.37";
    "E14" = "  -9.65%  ";
    "D15" = "6.588";
    "E15" = "  -7.14%  ";
    "D16" = "1.644.45";
    "E16" = "  -4.01%  ";
    "D17" = "0.00001037";
    "E17" = "  -8.94%  ";
    "D18" = "0.06564";
    "E18" = "  -2.81%  ";
    "D19" = "0.9997";
    "E19" = "  +0.21%  ";
    "D20" = "78.52";
    "E20" = "  -7.94%  ";
    "D21" = "5.926";
    "E21" = "  -7.86%  ";
    "D22" = "15.61";
    "E22" = "  -10.53%  ";
    "D23" = "12.53";
    "E23" = "  -5.77%  ";
    "D24" = "24.428.88";
    "E24" = "  -2.01%  ";
    "D25" = "2.456";
    "E25" = "  -0.02%  ";
    "D26" = "2.319";
    "E26" = "  -17.61%  ";
    "D27" = "146.35";
    "E27" = "  -3.36%  ";
    "D28" = "18.57";
    "E28" = "  -9.81%  ";
    "D29" = "1.827.14";
    "E29" = "  -3.99%  ";
    "B30" = "ImmutableX";
    "C30" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";
    "D30" = "1.188";
    "E30" = "  -5.29%  ";
    "B31" = "BitcoinCash";
    "C31" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch";
    "D31" = "124.00";
    "E31" = "  -7.08%  ";
    "D32" = "4.070";
    "E32" = "  -4.21%  ";
    "D33" = "5.668";
    "E33" = "  -18.37%  ";
    "D34" = "0.08413";
    "E34" = "  -5.27%  ";
    "D35" = "1.644";
    "E35" = "  -9.03%  ";
    "D36" = "12.08";
    "E36" = "  -13.81%  ";
    "B37" = "InternetComputer(DFINITY)";
    "C37" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";
    "D37" = "5.183";
    "E37" = "  -8.79%  ";
    "B38" = "TrustWalletToken";
    "C38" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";
    "D38" = "1.265";
    "E38" = "  -1.67%  ";
    "D39" = "0.06007";
    "E39" = "  -10.49%  ";
    "D40" = "0.02223";
    "E40" = "  -8.51%  ";
    "D41" = "0.2055";
    "E41" = "  -8.78%  ";
    "D42" = "8.098";
    "E42" = "  -13.98%  ";
    "D44" = "0.5883";
    "E44" = "  -9.52%  ";
    "D45" = "3.753";
    "E45" = "  -2.46%  ";
    "D46" = "12.49";
    "E46" = "  -11.23%  ";
    "D47" = "0.5603";
    "E47" = "  -9.69%  ";
    "D48" = "122.24";
    "E48" = "  -6.58%  ";
    "D49" = "1.942";
    "E49" = "  -9.72%  ";
    "D50" = "0.06904";
    "E50" = "  -5.80%  ";
    "D51" = "1.179";
    "E51" = "  -3.86%  ";
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    Set-TextValue $range $updates[$cellRef]
}

Write-Output "Applied $($updates.Count) cell updates"
